{"js": "// Apply the \"Added many more features\" edits to the D'Cirque review.\n// Each entry is an exact-text search-and-replace; matchCase keeps the\n// search from clobbering similar-but-different sentences elsewhere in\n// the document.\nconst replacements = [\n  [\n    \"Play D'Cirque Slot for Free - Dazzling Circus Theme\",\n    \"Play D'Cirque - Free Circus-Themed Slot\"\n  ],\n  [\n    \"Expanding symbols trigger respins and free spins\",\n    \"Expanding symbols and respins\"\n  ],\n  [\n    \"Visually stunning graphics and attention to detail\",\n    \"Free spins with wild symbols\"\n  ],\n  [\n    \"High RTP of 96.11%\",\n    \"Visually stunning graphics\"\n  ],\n  [\n    \"Potential winnings of up to 1,500 times your bet\",\n    \"Circus theme appeals to a diverse range of players\"\n  ],\n  [\n    \"Gameplay structure below par compared to other Peter & Sons slots\",\n    \"Gameplay structure is below par compared to other slots in the software's portfolio\"\n  ],\n  [\n    \"Limited number of paylines\",\n    \"Medium volatility may not appeal to all players\"\n  ],\n  [\n    \"D'Cirque by Peter & Sons is a circus-themed slot with expanding symbols and free spins. Try your luck with potential winnings of up to 1,500 times your bet.\",\n    \"Discover the visually stunning D'Cirque slot game and play for free. Win big with expanding symbols and free spins.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the D'Cirque review.\n# Each entry is an exact-text search-and-replace. We locate each hit with\n# Find.Execute (plain text, no wildcards) and then assign the replacement\n# directly to the found Range's .Text \u2014 going through Find's own\n# Replacement.Text parameter triggers Word's AutoFormat \"smart quotes\"\n# substitution (straight ' becomes a curly U+2019), which the target text\n# must NOT have, so we avoid that code path entirely. The outer while loop\n# repeats the find/replace until no more matches remain, so occurrences\n# that appear more than once (the heading) all get updated.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play D'Cirque Slot for Free - Dazzling Circus Theme\", \"Play D'Cirque - Free Circus-Themed Slot\"),\n    @(\"Expanding symbols trigger respins and free spins\", \"Expanding symbols and respins\"),\n    @(\"Visually stunning graphics and attention to detail\", \"Free spins with wild symbols\"),\n    @(\"High RTP of 96.11%\", \"Visually stunning graphics\"),\n    @(\"Potential winnings of up to 1,500 times your bet\", \"Circus theme appeals to a diverse range of players\"),\n    @(\"Gameplay structure below par compared to other Peter & Sons slots\", \"Gameplay structure is below par compared to other slots in the software's portfolio\"),\n    @(\"Limited number of paylines\", \"Medium volatility may not appeal to all players\"),\n    @(\"D'Cirque by Peter & Sons is a circus-themed slot with expanding symbols and free spins. Try your luck with potential winnings of up to 1,500 times your bet.\", \"Discover the visually stunning D'Cirque slot game and play for free. Win big with expanding symbols and free spins.\")\n)\n\nforeach ($pair in $replacements) {\n    $target = $pair[0]\n    $replacement = $pair[1]\n\n    $continue = $true\n    while ($continue) {\n        $rng = $d.Content\n        $find = $rng.Find\n        $find.Text = $target\n        $find.MatchCase = $true\n        $find.MatchWildcards = $false\n        $found = $find.Execute()\n        if ($found) {\n            $rng.Text = $replacement\n        } else {\n            $continue = $false\n        }\n    }\n}\n\nWrite-Output \"done\"\n"}
